$wb = $excel.ActiveWorkbook


# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3580
$ws1.Range("F4").Value = 380
$ws1.Range("F5").Value = 8353
$ws1.Range("F7").Value = 131
$ws1.Range("F8").Value = 2244
$ws1.Range("F10").Value = 107
$ws1.Range("C11").Value = "上海·原神X星穹铁道x绝区零同人ONLY（取消）"
$ws1.Range("D11").Value = "顾村镇蕰川路6号 智慧湾科创园"
$ws1.Range("E11").Value = "2024.10.01 10:00-10.02 17:00"
$ws1.Range("F11").Value = 510
$ws1.Range("G11").Value = "不可售"
$ws1.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=90135"
$ws1.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202407/FF8HGnt01722418798545.jpeg"
$ws1.Range("C12").Value = "上海·国潮新次元——2024南翔国潮大会（免费）"
$ws1.Range("D12").Value = "南翔镇金迈路槎溪路路口 银翔湖公园"
$ws1.Range("E12").Value = "2024.10.01 10:00-10.02 20:00"
$ws1.Range("F12").Value = 79
$ws1.Range("G12").Value = 20
$ws1.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=92824"
$ws1.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202409/WRi6DW941727096142734.jpeg"
$ws1.Range("C13").Value = "上海·妖漫第六十届动漫展（免费漫展）"
$ws1.Range("D13").Value = "秀浦路668号 新田360广场(上海康桥店)"
$ws1.Range("E13").Value = "2024.10.01 12:00-10.03 18:00"
$ws1.Range("F13").Value = 667
$ws1.Range("G13").Value = 35.8
$ws1.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=92641"
$ws1.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202409/QXoONdAk1726734646169.png"
$ws1.Range("C14").Value = "上海·崖牙首次个人画展 - 金色生灵GOLDEN CREATURE"
$ws1.Range("D14").Value = "福州路390号 外文书店"
$ws1.Range("E14").Value = "2024.10.01 11:00-10.14 17:00"
$ws1.Range("F14").Value = 112
$ws1.Range("G14").Value = 10
$ws1.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=92513"
$ws1.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202409/BKsTihwt1726310450103.png"
$ws1.Range("C15").Value = "上海·异星宇宙·综合同人Only动漫展"
$ws1.Range("D15").Value = "中山北路3300号 上海环球港"
$ws1.Range("E15").Value = "2024.10.01 09:30-10.05 17:00"
$ws1.Range("F15").Value = 7527
$ws1.Range("G15").Value = 60
$ws1.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=92752"
$ws1.Range("I15").Value = "//i2.hdslb.com/bfs/openplatform/202409/gKL5QFbk1727085314378.jpeg"
$ws1.Range("F16").Value = 7741
$ws1.Range("F18").Value = 57967
$ws1.Range("F19").Value = 57967
$ws1.Range("F20").Value = 4861
$ws1.Range("F21").Value = 1064
$ws1.Range("F22").Value = 965
$ws1.Range("F23").Value = 519
$ws1.Range("F26").Value = 15
$ws1.Range("F27").Value = 620
$ws1.Range("F28").Value = 5322
$ws1.Range("F29").Value = 603
$ws1.Range("F31").Value = 55
$ws1.Range("F32").Value = 927
$ws1.Range("F33").Value = 1420
$ws1.Range("F34").Value = 2016
$ws1.Range("F37").Value = 238
$ws1.Range("F39").Value = 6
$ws1.Range("F43").Value = 543
$ws1.Range("F44").Value = 271
$ws1.Range("F48").Value = 17

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 205
$ws2.Range("F3").Value = 39
$ws2.Range("F5").Value = 64
$ws2.Range("F6").Value = 161
$ws2.Range("F9").Value = 7673
$ws2.Range("F14").Value = 9
$ws2.Range("F27").Value = 133
$ws2.Range("F45").Value = 43

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 2387
$ws3.Range("F5").Value = 1623
$ws3.Range("F7").Value = 685
$ws3.Range("F8").Value = 2430
$ws3.Range("F9").Value = 9453
$ws3.Range("F10").Value = 1778
$ws3.Range("F15").Value = 281
$ws3.Range("F16").Value = 2447
$ws3.Range("F17").Value = 192
$ws3.Range("F18").Value = 70
$ws3.Range("F19").Value = 538

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3580
$ws4.Range("F4").Value = 2387
$ws4.Range("F5").Value = 8353
$ws4.Range("F6").Value = 1778
$ws4.Range("F8").Value = 281
$ws4.Range("F9").Value = 192
$ws4.Range("F10").Value = 79
$ws4.Range("F11").Value = 667
$ws4.Range("F12").Value = 7741
$ws4.Range("F13").Value = 57967
$ws4.Range("F14").Value = 205
$ws4.Range("F15").Value = 39
$ws4.Range("F16").Value = 4861
$ws4.Range("F17").Value = 965
$ws4.Range("F18").Value = 519
$ws4.Range("F20").Value = 620
$ws4.Range("F21").Value = 162
$ws4.Range("F23").Value = 927
$ws4.Range("F24").Value = 1420
$ws4.Range("F25").Value = 2016
$ws4.Range("F27").Value = 538
$ws4.Range("F29").Value = 9
$ws4.Range("F32").Value = 238
$ws4.Range("F36").Value = 543
$ws4.Range("F49").Value = 43

Write-Output "Done applying updates"
